$d = $word.ActiveDocument

# The target table row (Table 2, Row 2) originally contains, in the left
# cell, five paragraphs listing individual names (Abner, Isabella, Jose,
# Laura, Vinicius) and in the right cell a one-paragraph description.
# The edit collapses the five name paragraphs down to a single paragraph
# reading "Desenvolvedores", and rewords the description text.

$w_ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the five paragraphs by their known text, searching the document
# paragraphs collection (works reliably right after the document is
# freshly loaded, before any structural edits shift indices).
$namesText = @(
    "Abner de Melo Porto;",
    "Isabella Mota Franco;",
    "José Alves de Oliveira;",
    "Laura Tazue Tavares Chirazawa;",
    "Vinicius Tertuliano da Silva."
)

$startIndex = -1
for ($p = 1; $p -le $d.Paragraphs.Count; $p++) {
    $para = $d.Paragraphs.Item($p)
    if ($para.Range.Text -eq ($namesText[0] + "`r")) {
        $startIndex = $p
        break
    }
}

if ($startIndex -eq -1) {
    throw "Could not locate the 'Abner de Melo Porto;' paragraph"
}

# Delete the first four name paragraphs (Abner, Isabella, Jose, Laura),
# working from the last one backwards so earlier indices stay valid.
for ($i = 3; $i -ge 0; $i--) {
    $para = $d.Paragraphs.Item($startIndex + $i)
    $para.Range.Delete()
}

# The "Vinicius Tertuliano da Silva." paragraph is now at $startIndex;
# replace its whole contents (including the paragraph mark) with a
# single paragraph reading "Desenvolvedores", preserving the paragraph
# formatting (pPr) and using a simple (no explicit size) run format.
$paraVinicius = $d.Paragraphs.Item($startIndex)
$rngVinicius = $d.Range($paraVinicius.Range.Start, $paraVinicius.Range.End)
$xmlDesenvolvedores = "<w:p $w_ns><w:pPr><w:widowControl w:val=`"0`"/><w:spacing w:line=`"240`" w:lineRule=`"auto`"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr><w:t xml:space=`"preserve`">Desenvolvedores</w:t></w:r></w:p>"
$rngVinicius.InsertXML($xmlDesenvolvedores)

# The description paragraph ("Desenvolvedores do Website e do App.") is
# now the very next paragraph; replace its whole contents with the new
# wording, preserving the paragraph formatting and the trailing empty
# run that followed the text run in the original document.
$paraDesc = $d.Paragraphs.Item($startIndex + 1)
$rngDesc = $d.Range($paraDesc.Range.Start, $paraDesc.Range.End)
$xmlDesc = "<w:p $w_ns><w:pPr><w:widowControl w:val=`"0`"/><w:spacing w:line=`"240`" w:lineRule=`"auto`"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr><w:t xml:space=`"preserve`">Equipe que irá desenvolver o Website e App Softlife.</w:t></w:r><w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr></w:r></w:p>"
$rngDesc.InsertXML($xmlDesc)
